# ---------------------------------------------------------------------------
# PropostaSIColetaLixo.docx edits
#   1. Title paragraph: merge the two "Projeto " / "SIColetaLixo" runs into a
#      single run (drops the spell-check proofErr markers along the way).
#   2. "1a etapa" paragraph: rework the sentence about software adjustments.
#   3. Append a new "Metodo Agil" section at the end of the document.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title: "Projeto " + "SIColetaLixo" (two runs, proofErr-wrapped) -> one
#    clean run reading "Projeto SIColetaLixo".
# ---------------------------------------------------------------------------
$titleRange = $d.Range(0, 21)
$titleRange.Delete()

$d.Range(0, 0).InsertParagraphBefore()

$newTitle = $d.Range(0, 0)
$newTitle.InsertBefore("Projeto SIColetaLixo")

$titleRun = $d.Range(0, 20)
$titleRun.Font.Bold = $true
$titleRun.Font.Size = 16

# ---------------------------------------------------------------------------
# 2) "1a etapa" paragraph rewrite.
# ---------------------------------------------------------------------------
$oldSentence = "no software para melhor se integrar as suas atividades."
$newSentence = "no escopo do software e adição de recursos que melhor se integrem as suas atividades."
$d.Content.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, $newSentence, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Append "Método Ágil" section after the existing trailing blank
#    paragraphs.
# ---------------------------------------------------------------------------
$lastParaIndex = $d.Paragraphs.Count
$d.Paragraphs($lastParaIndex).Range.InsertParagraphAfter()

$headingPara = $d.Paragraphs($lastParaIndex + 1).Range
$headingPara.InsertBefore("Método Ágil:")
$headingPara.Font.Bold = $true
$headingPara.Font.Size = 14

$headingPara.InsertParagraphAfter()
$para1 = $d.Paragraphs($lastParaIndex + 2).Range
$para1.InsertBefore("Como o cliente quer entregas rápidas que o permitem ter acesso o quanto antes a ferramentas do sistema, a utilização de um método ágil de produção como o SCRUM é recomendada.")
$para1.Font.Bold = $false
$para1.Font.Size = 12

$para1.InsertParagraphAfter()
$para2 = $d.Paragraphs($lastParaIndex + 3).Range
$para2.InsertBefore("Com ele, diversos times podem ser feitos focando em diferentes recursos necessários para o software. Com as entregas e o feedback do cliente, novos Sprints podem ser feitos por cada time para melhorar as funções do sistema e adicionar novas ferramentas, além de uma nova análise e constante atualização do escopo se necessário.")
$para2.Font.Bold = $false
$para2.Font.Size = 12

Write-Output "done"
